$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032252136799856
$ws.Range("D2").Value = 1.033248583514515
$ws.Range("E2").Value = 1.039769076285848
$ws.Range("F2").Value = 1.046118930338881
$ws.Range("I2").Value = 1.025700725390273
$ws.Range("J2").Value = 1.037383380229097
$ws.Range("K2").Value = 1.036051478509139
$ws.Range("L2").Value = 1.042553323767552
$ws.Range("M2").Value = 1.04888526071227
$ws.Range("N2").Value = 1.016376239474715
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03426545083285
$ws.Range("D3").Value = 1.035056343937699
$ws.Range("E3").Value = 1.041645668744029
$ws.Range("F3").Value = 1.048312227379284
$ws.Range("I3").Value = 1.025940528533644
$ws.Range("J3").Value = 1.039033510198414
$ws.Range("K3").Value = 1.03766551470421
$ws.Range("L3").Value = 1.04423739605583
$ws.Range("M3").Value = 1.05088654771063
$ws.Range("N3").Value = 1.016923269456327
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.035555253493322
$ws.Range("D4").Value = 1.036214436524599
$ws.Range("E4").Value = 1.042845136549768
$ws.Range("F4").Value = 1.049707704106653
$ws.Range("I4").Value = 1.026084975133559
$ws.Range("J4").Value = 1.040088242263495
$ws.Range("K4").Value = 1.038697859055451
$ws.Range("L4").Value = 1.045311912964725
$ws.Range("M4").Value = 1.052157490626605
$ws.Range("N4").Value = 1.017272843634881
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.036094442021017
$ws.Range("D5").Value = 1.036698558276703
$ws.Range("E5").Value = 1.043345901714556
$ws.Range("F5").Value = 1.050288757205262
$ws.Range("I5").Value = 1.026143146794347
$ws.Range("J5").Value = 1.040528584891023
$ws.Range("K5").Value = 1.039129019223672
$ws.Range("L5").Value = 1.045760055748595
$ws.Range("M5").Value = 1.052686116223936
$ws.Range("N5").Value = 1.01741876960325
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.036184797167301
$ws.Range("D6").Value = 1.036779685123694
$ws.Range("E6").Value = 1.0434297792262
$ws.Range("F6").Value = 1.05038599202963
$ws.Range("I6").Value = 1.026152764683142
$ws.Range("J6").Value = 1.040602341889456
$ws.Range("K6").Value = 1.039201247832036
$ws.Range("L6").Value = 1.045835092199506
$ws.Range("M6").Value = 1.05277454381731
$ws.Range("N6").Value = 1.017443210977894
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.035562470047767
$ws.Range("D7").Value = 1.036220916087356
$ws.Range("E7").Value = 1.042851841443845
$ws.Range("F7").Value = 1.049715490091868
$ws.Range("I7").Value = 1.026085762443314
$ws.Range("J7").Value = 1.040094138129943
$ws.Range("K7").Value = 1.038703631327227
$ws.Range("L7").Value = 1.045317915085016
$ws.Range("M7").Value = 1.052164576366942
$ws.Range("N7").Value = 1.017274797550558
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032935260034086
$ws.Range("D8").Value = 1.03386196654224
$ws.Range("E8").Value = 1.040406380646978
$ws.Range("F8").Value = 1.046865126697598
$ws.Range("I8").Value = 1.025783996397267
$ws.Range("J8").Value = 1.037943773548519
$ws.Range("K8").Value = 1.036599472067835
$ws.Range("L8").Value = 1.043125641580771
$ws.Range("M8").Value = 1.049566624124195
$ws.Range("N8").Value = 1.016562029915166
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028204050616421
$ws.Range("D9").Value = 1.029613693176449
$ws.Range("E9").Value = 1.035981110292742
$ws.Range("F9").Value = 1.041657024631261
$ws.Range("I9").Value = 1.025169469285779
$ws.Range("J9").Value = 1.034052610371359
$ws.Range("K9").Value = 1.032797245622586
$ws.Range("L9").Value = 1.039143749692325
$ws.Range("M9").Value = 1.044801258229978
$ws.Range("N9").Value = 1.015271652466504
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024977738195995
$ws.Range("D10").Value = 1.026716640656763
$ws.Range("E10").Value = 1.032949065558033
$ws.Range("F10").Value = 1.038054937020859
$ws.Range("I10").Value = 1.024703110450822
$ws.Range("J10").Value = 1.031386596988061
$ws.Range("K10").Value = 1.030195738813927
$ws.Range("L10").Value = 1.036405600294665
$ws.Range("M10").Value = 1.041493209578412
$ws.Range("N10").Value = 1.014387160495711
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023562700111043
$ws.Range("D11").Value = 1.025446007154784
$ws.Range("E11").Value = 1.031615829815181
$ws.Range("F11").Value = 1.036463073555343
$ws.Range("I11").Value = 1.024487473943922
$ws.Range("J11").Value = 1.030214332395829
$ws.Range("K11").Value = 1.029052682203049
$ws.Range("L11").Value = 1.035199262676302
$ws.Range("M11").Value = 1.040028461268431
$ws.Range("N11").Value = 1.013998151440986
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023034304037809
$ws.Range("D12").Value = 1.02497153359019
$ws.Range("E12").Value = 1.031117468662883
$ws.Range("F12").Value = 1.035866840443401
$ws.Range("I12").Value = 1.024405294651763
$ws.Range("J12").Value = 1.029776145536082
$ws.Range("K12").Value = 1.028625539776245
$ws.Range("L12").Value = 1.034747986393466
$ws.Range("M12").Value = 1.039479418589331
$ws.Range("N12").Value = 1.0138527281223
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023147774275446
$ws.Range("D13").Value = 1.02507342425928
$ws.Range("E13").Value = 1.031224512187213
$ws.Range("F13").Value = 1.035994960132801
$ws.Range("I13").Value = 1.024423017021295
$ws.Range("J13").Value = 1.029870263993233
$ws.Range("K13").Value = 1.028717280268619
$ws.Range("L13").Value = 1.034844932293791
$ws.Range("M13").Value = 1.039597416854367
$ws.Range("N13").Value = 1.013883964309158
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023519080076398
$ws.Range("D14").Value = 1.025406838521111
$ws.Range("E14").Value = 1.031574699631503
$ws.Range("F14").Value = 1.036413890292804
$ws.Range("I14").Value = 1.024480723601231
$ws.Range("J14").Value = 1.030178168374869
$ws.Range("K14").Value = 1.029017427113265
$ws.Range("L14").Value = 1.035162025557674
$ws.Range("M14").Value = 1.039983179277948
$ws.Range("N14").Value = 1.013986149777962
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023747481949644
$ws.Range("D15").Value = 1.025611932098667
$ws.Range("E15").Value = 1.03179004332908
$ws.Range("F15").Value = 1.036671348269388
$ws.Range("I15").Value = 1.0245160018584
$ws.Range("J15").Value = 1.030367511075592
$ws.Range("K15").Value = 1.029202016126633
$ws.Range("L15").Value = 1.035356972218185
$ws.Range("M15").Value = 1.04022019825445
$ws.Range("N15").Value = 1.014048985910909
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025071263124915
$ws.Range("D16").Value = 1.026800621293001
$ws.Range("E16").Value = 1.033037112555189
$ws.Range("F16").Value = 1.038159896916519
$ws.Range("I16").Value = 1.024717130880048
$ws.Range("J16").Value = 1.031464013915815
$ws.Range("K16").Value = 1.030271244491141
$ws.Range("L16").Value = 1.036485218175436
$ws.Range("M16").Value = 1.041589729358963
$ws.Range("N16").Value = 1.01441284893718
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025896755271934
$ws.Range("D17").Value = 1.027541870659687
$ws.Range("E17").Value = 1.03381386369938
$ws.Range("F17").Value = 1.039084940380523
$ws.Range("I17").Value = 1.024839609486917
$ws.Range("J17").Value = 1.032146989035638
$ws.Range("K17").Value = 1.030937455176434
$ws.Range("L17").Value = 1.037187340489521
$ws.Range("M17").Value = 1.042440064757872
$ws.Range("N17").Value = 1.014639462859229
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026376517469418
$ws.Range("D18").Value = 1.02797267169874
$ws.Range("E18").Value = 1.034264972706947
$ws.Range("F18").Value = 1.039621410105515
$ws.Range("I18").Value = 1.024909729045779
$ws.Range("J18").Value = 1.032543638362052
$ws.Range("K18").Value = 1.03132444897694
$ws.Range("L18").Value = 1.037594884745109
$ws.Range("M18").Value = 1.042932938337951
$ws.Range("N18").Value = 1.014771063929198
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026539812336336
$ws.Range("D19").Value = 1.028119301704684
$ws.Range("E19").Value = 1.034418459742829
$ws.Range("F19").Value = 1.03980381130804
$ws.Range("I19").Value = 1.024933414785606
$ws.Range("J19").Value = 1.032678596168593
$ws.Range("K19").Value = 1.031456135237507
$ws.Range("L19").Value = 1.03773351128459
$ws.Range("M19").Value = 1.043100470803958
$ws.Range("N19").Value = 1.014815838991409
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025808367638961
$ws.Range("D20").Value = 1.027462503188586
$ws.Range("E20").Value = 1.033730728634468
$ws.Range("F20").Value = 1.038986012661765
$ws.Range("I20").Value = 1.024826605408214
$ws.Range("J20").Value = 1.032073890530873
$ws.Range("K20").Value = 1.030866142601187
$ws.Range("L20").Value = 1.037112215989463
$ws.Range("M20").Value = 1.042349154559873
$ws.Range("N20").Value = 1.014615209387478
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023409817393968
$ws.Range("D21").Value = 1.025308726033996
$ws.Range("E21").Value = 1.031471665461338
$ws.Range("F21").Value = 1.036290663246122
$ws.Range("I21").Value = 1.024463788134388
$ws.Range("J21").Value = 1.030087574870419
$ws.Range("K21").Value = 1.02892911258971
$ws.Range("L21").Value = 1.035068738171669
$ws.Range("M21").Value = 1.039869719973336
$ws.Range("N21").Value = 1.013956084519799
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021885581266746
$ws.Range("D22").Value = 1.023940037490513
$ws.Range("E22").Value = 1.030033103529424
$ws.Range("F22").Value = 1.034567323000504
$ws.Range("I22").Value = 1.0242236098894
$ws.Range("J22").Value = 1.028822718284323
$ws.Range("K22").Value = 1.027696375626781
$ws.Range("L22").Value = 1.033765431481132
$ws.Range("M22").Value = 1.038281987341769
$ws.Range("N22").Value = 1.013536284590724
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022695168332372
$ws.Range("D23").Value = 1.024667006484618
$ws.Range("E23").Value = 1.030797465964464
$ws.Range("F23").Value = 1.03548365589721
$ws.Range("I23").Value = 1.024352084591047
$ws.Range("J23").Value = 1.029494782401551
$ws.Range("K23").Value = 1.028351304006797
$ws.Range("L23").Value = 1.034458119049775
$ws.Range("M23").Value = 1.039126444230045
$ws.Range("N23").Value = 1.013759346888431
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025848311557466
$ws.Range("D24").Value = 1.027498370740167
$ws.Range("E24").Value = 1.033768299833666
$ws.Range("F24").Value = 1.039030723388641
$ws.Range("I24").Value = 1.02483248546957
$ws.Range("J24").Value = 1.032106925904836
$ws.Range("K24").Value = 1.030898370610973
$ws.Range("L24").Value = 1.037146167667583
$ws.Range("M24").Value = 1.042390242576918
$ws.Range("N24").Value = 1.014626170275145
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029439614656974
$ws.Range("D25").Value = 1.030723151865712
$ws.Range("E25").Value = 1.037139284585103
$ws.Range("F25").Value = 1.043025948746921
$ws.Range("I25").Value = 1.025338248163195
$ws.Range("J25").Value = 1.035070986013849
$ws.Range("K25").Value = 1.033791722185767
$ws.Range("L25").Value = 1.040187610817905
$ws.Range("M25").Value = 1.046055940755048
$ws.Range("N25").Value = 1.015609433032043
